$d = $word.ActiveDocument

# --- Paragraph 1: title (no text change, just style normalization) ---
$p1 = $d.Paragraphs(1)
$p1.Style = "Normal"
$p1.Format.LineSpacingRule = 5
$p1.Format.LineSpacing = 24
$p1.Alignment = 1

# --- Paragraph 2: Kaden Roof -- collapse the several runs into one ---
$d.Content.Find.Execute(
    "Kaden Roof: Created/worked on Phase 2 Presentation, worked on the Django files. Created this file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kaden Roof: Created/worked on Phase 2 Presentation, worked on the Django files. Created this file.",
    2)
$p2 = $d.Paragraphs(2)
$p2.Style = "Normal"
$p2.Format.LineSpacingRule = 5
$p2.Format.LineSpacing = 24

# --- Paragraph 3: Grant Burkemper -- collapse the two runs into one ---
$d.Content.Find.Execute(
    "Grant Burkemper: Created database, worked on Gantt documentation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Grant Burkemper: Created database, worked on Gantt documentation",
    2)
$p3 = $d.Paragraphs(3)
$p3.Style = "Normal"
$p3.Format.LineSpacingRule = 5
$p3.Format.LineSpacing = 24

# --- Paragraph 4: previously empty, fill in Jacob Nagel's contribution ---
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "Jacob Nagel (primarily Server/Backend): Initially began server/backend research, created server notes document listing basic/preliminary info and SQL commands that may be used, wrote most of system design document, completely reworked and wrote most of milestone list.  Ready to begin coding and implementing server and back end in phase 3 this week."
$p4b = $d.Paragraphs.Last
$p4b.Style = "Normal"
$p4b.Format.LineSpacingRule = 5
$p4b.Format.LineSpacing = 24
$p4b.Format.SpaceBefore = 0
$p4b.Format.SpaceAfter = 8
